$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.049.72'
$ws.Range("E2").Value = '  +2.80%  '
$ws.Range("D3").Value = '3.794.73'
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = "'696.88"
$ws.Range("E5").Value = '  +10.52%  '
$ws.Range("D6").Value = "'172.75"
$ws.Range("E6").Value = '  +3.92%  '
$ws.Range("D7").Value = '3.794.31'
$ws.Range("E7").Value = '  +0.76%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E10").Value = '  +2.71%  '
$ws.Range("D11").Value = "'7.46"
$ws.Range("E11").Value = '  +10.92%  '
$ws.Range("D12").Value = "'0.461"
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("E13").Value = '  +7.77%  '
$ws.Range("D14").Value = "'36.16"
$ws.Range("E14").Value = '  +3.62%  '
$ws.Range("D15").Value = '4.438.63'
$ws.Range("E15").Value = '  +0.71%  '
$ws.Range("D16").Value = '3.814.18'
$ws.Range("E16").Value = '  +1.12%  '
$ws.Range("D17").Value = '71.116.58'
$ws.Range("D18").Value = "'17.84"
$ws.Range("E18").Value = '  +1.60%  '
$ws.Range("D19").Value = "'7.20"
$ws.Range("E19").Value = '  +2.99%  '
$ws.Range("D20").Value = "'0.114"
$ws.Range("E20").Value = '  +1.09%  '
$ws.Range("D21").Value = "'11.13"
$ws.Range("E21").Value = '  +17.11%  '
$ws.Range("D22").Value = "'482.59"
$ws.Range("E22").Value = '  +4.35%  '
$ws.Range("D23").Value = "'0.715"
$ws.Range("E23").Value = '  +1.49%  '
$ws.Range("D24").Value = "'83.82"
$ws.Range("E24").Value = '  +2.19%  '
$ws.Range("D25").Value = "'0.0000145"
$ws.Range("E25").Value = '  +0.67%  '
$ws.Range("D26").Value = "'12.36"
$ws.Range("E26").Value = '  +2.47%  '
$ws.Range("D27").Value = "'10.52"
$ws.Range("E27").Value = '  +4.72%  '
$ws.Range("D28").Value = "'2.18"
$ws.Range("E28").Value = '  +1.87%  '
$ws.Range("D29").Value = '3.950.40'
$ws.Range("E29").Value = '  +0.74%  '
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("E31").Value = '  +14.42%  '
$ws.Range("E32").Value = '  +0.26%  '
$ws.Range("D33").Value = "'7.56"
$ws.Range("E33").Value = '  +7.06%  '
$ws.Range("D34").Value = "'29.56"
$ws.Range("E34").Value = '  +3.80%  '
$ws.Range("E35").Value = '  -1.21%  '
$ws.Range("D36").Value = "'9.20"
$ws.Range("E36").Value = '  +3.02%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").Value = '3.750.37'
$ws.Range("E38").Value = '  +0.70%  '
$ws.Range("E39").Value = '  +2.18%  '
$ws.Range("D40").Value = "'3.48"
$ws.Range("E40").Value = '  +5.91%  '
$ws.Range("E41").Value = '  +3.28%  '
$ws.Range("D42").Value = "'2.22"
$ws.Range("E42").Value = '  +11.92%  '
$ws.Range("E43").Value = '  +22.99%  '
$ws.Range("D44").Value = "'0.964"
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").Value = "'45.49"
$ws.Range("E47").Value = '  +5.10%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = "'160.65"
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").Value = "'48.96"
$ws.Range("E49").Value = '  +4.09%  '
$ws.Range("E50").Value = '  -1.53%  '
$ws.Range("D51").Value = "'0.299"
$ws.Range("E51").Value = '  +1.44%  '
